# Auto-generated Excel COM-interop script to apply the NEW_JERSEY_2024-style cleanup edit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header columns to snake_case English names
$ws.Range("A1").Value = 'mx_state'
$ws.Range("B1").Value = 'mx_municipality'
$ws.Range("C1").Value = 'n_matriculas'
$ws.Range("D1").Value = 'pct_matriculas'

# 2. Title-case the Spanish connector words ("de", "del", "la", ...) in state/municipality names
$ws.Range("B3").Value = 'Rincón De Romos'
$ws.Range("B23").Value = 'Amatenango De La Frontera'
$ws.Range("B26").Value = 'Bejucal De Ocampo'
$ws.Range("B28").Value = 'Benemérito De Las Américas'
$ws.Range("B34").Value = 'Chiapa De Corzo'
$ws.Range("B38").Value = 'Comitán De Domínguez'
$ws.Range("B59").Value = 'Mazapa De Madero'
$ws.Range("B64").Value = 'Ocozocoautla De Espinosa'
$ws.Range("B70").Value = 'San Cristóbal De Las Casas'
$ws.Range("A102").Value = 'Ciudad De México'
$ws.Range("B105").Value = 'Cuajimalpa De Morelos'
$ws.Range("A120").Value = 'Coahuila De Zaragoza'
$ws.Range("B124").Value = 'San Juan De Sabinas'
$ws.Range("B134").Value = 'Nombre De Dios'
$ws.Range("A139").Value = 'Estado De México'
$ws.Range("B139").Value = 'Acambay De Ruíz Castañeda'
$ws.Range("B142").Value = 'Almoloya De Alquisiras'
$ws.Range("B143").Value = 'Almoloya De Juárez'
$ws.Range("B148").Value = 'Atizapán De Zaragoza'
$ws.Range("B154").Value = 'Chapa De Mota'
$ws.Range("B158").Value = 'Coacalco De Berriozábal'
$ws.Range("B165").Value = 'Ecatepec De Morelos'
$ws.Range("B171").Value = 'Ixtapan De La Sal'
$ws.Range("B172").Value = 'Ixtapan Del Oro'
$ws.Range("B182").Value = 'Naucalpan De Juárez'
$ws.Range("B189").Value = 'San Felipe Del Progreso'
$ws.Range("B190").Value = 'San José Del Rincón'
$ws.Range("B191").Value = 'San Martín De Las Pirámides'
$ws.Range("B193").Value = 'San Simón De Guerrero'
$ws.Range("B203").Value = 'Tenango Del Valle'
$ws.Range("B212").Value = 'Tlalnepantla De Baz'
$ws.Range("B217").Value = 'Valle De Bravo'
$ws.Range("B218").Value = 'Valle De Chalco Solidaridad'
$ws.Range("B221").Value = 'Villa De Allende'
$ws.Range("B222").Value = 'Villa Del Carbón'
$ws.Range("B230").Value = 'Apaseo El Alto'
$ws.Range("B236").Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range("B246").Value = 'San Francisco Del Rincón'
$ws.Range("B247").Value = 'San Luis De La Paz'
$ws.Range("B248").Value = 'San Miguel De Allende'
$ws.Range("B249").Value = 'Santa Cruz De Juventino Rosas'
$ws.Range("B250").Value = 'Silao De La Victoria'
$ws.Range("B254").Value = 'Valle De Santiago'
$ws.Range("B259").Value = 'Acapulco De Juárez'
$ws.Range("B262").Value = 'Ajuchitlán Del Progreso'
$ws.Range("B263").Value = 'Alcozauca De Guerrero'
$ws.Range("B266").Value = 'Atenango Del Río'
$ws.Range("B267").Value = 'Atlamajalcingo Del Monte'
$ws.Range("B269").Value = 'Atoyac De Álvarez'
$ws.Range("B270").Value = 'Ayutla De Los Libres'
$ws.Range("B273").Value = 'Chilapa De Álvarez'
$ws.Range("B274").Value = 'Chilpancingo De Los Bravo'
$ws.Range("B279").Value = 'Coyuca De Benítez'
$ws.Range("B283").Value = 'Cuetzala Del Progreso'
$ws.Range("B289").Value = 'Huitzuco De Los Figueroa'
$ws.Range("B290").Value = 'Iguala De La Independencia'
$ws.Range("B292").Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range("B298").Value = 'Mártir De Cuilapan'
$ws.Range("B308").Value = 'Taxco De Alarcón'
$ws.Range("B311").Value = 'Tepecoacuilco De Trujano'
$ws.Range("B312").Value = 'Tixtla De Guerrero'
$ws.Range("B315").Value = 'Tlalixtaquilla De Maldonado'
$ws.Range("B316").Value = 'Tlapa De Comonfort'
$ws.Range("B318").Value = 'Técpan De Galeana'
$ws.Range("B323").Value = 'Zihuatanejo De Azueta'
$ws.Range("B333").Value = 'Atotonilco De Tula'
$ws.Range("B334").Value = 'Atotonilco El Grande'
$ws.Range("B339").Value = 'Cuautepec De Hinojosa'
$ws.Range("B342").Value = 'Huasca De Ocampo'
$ws.Range("B346").Value = 'Huejutla De Reyes'
$ws.Range("B354").Value = 'Mineral Del Chico'
$ws.Range("B355").Value = 'Mixquiahuala De Juárez'
$ws.Range("B356").Value = 'Molango De Escamilla'
$ws.Range("B357").Value = 'Nopala De Villagrán'
$ws.Range("B358").Value = 'Omitlán De Juárez'
$ws.Range("B359").Value = 'Pachuca De Soto'
$ws.Range("B362").Value = 'Progreso De Obregón'
$ws.Range("B365").Value = 'Santiago Tulantepec De Lugo Guerrero'
$ws.Range("B366").Value = 'Santiago De Anaya'
$ws.Range("B369").Value = 'Tenango De Doria'
$ws.Range("B371").Value = 'Tepehuacán De Guerrero'
$ws.Range("B372").Value = 'Tepeji Del Río De Ocampo'
$ws.Range("B374").Value = 'Tezontepec De Aldama'
$ws.Range("B380").Value = 'Tula De Allende'
$ws.Range("B381").Value = 'Tulancingo De Bravo'
$ws.Range("B383").Value = 'Zacualtipán De Ángeles'
$ws.Range("B386").Value = 'Autlán De Navarro'
$ws.Range("B392").Value = 'Encarnación De Díaz'
$ws.Range("B396").Value = 'Lagos De Moreno'
$ws.Range("B402").Value = 'San Juan De Los Lagos'
$ws.Range("B405").Value = 'Tamazula De Gordiano'
$ws.Range("B410").Value = 'Zacoalco De Torres'
$ws.Range("B412").Value = 'Zapotlán El Grande'
$ws.Range("A414").Value = 'Michoacán De Ocampo'
$ws.Range("B424").Value = 'Coalcomán De Vázquez Pallares'
$ws.Range("B461").Value = 'Tiquicheo De Nicolás Romero'
$ws.Range("B490").Value = 'Puente De Ixtla'
$ws.Range("B494").Value = 'Tetela Del Volcán'
$ws.Range("B496").Value = 'Tlaltizapán De Zapata'
$ws.Range("B505").Value = 'Zacualpan De Amilpas'
$ws.Range("B507").Value = 'Bahía De Banderas'
$ws.Range("B527").Value = 'Acatlán De Pérez Figueroa'
$ws.Range("B532").Value = 'Capulálpam De Méndez'
$ws.Range("B533").Value = 'Chalcatongo De Hidalgo'
$ws.Range("B535").Value = 'Constancia Del Rosario'
$ws.Range("B537").Value = 'Cuilápam De Guerrero'
$ws.Range("B539").Value = 'El Barrio De La Soledad'
$ws.Range("B540").Value = 'Fresnillo De Trujano'
$ws.Range("B541").Value = 'Guevea De Humboldt'
$ws.Range("B542").Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range("B543").Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range("B544").Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range("B545").Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range("B546").Value = 'Huajuapan De León'
$ws.Range("B548").Value = 'Huautla De Jiménez'
$ws.Range("B549").Value = 'Ixtlán De Juárez'
$ws.Range("B557").Value = 'Mariscala De Juárez'
$ws.Range("B559").Value = 'Mazatlán Villa De Flores'
$ws.Range("B561").Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range("B562").Value = 'Mixistlán De La Reforma'
$ws.Range("B564").Value = 'Mártires De Tacubaya'
$ws.Range("B565").Value = 'Nejapa De Madero'
$ws.Range("B567").Value = 'Oaxaca De Juárez'
$ws.Range("B568").Value = 'Ocotlán De Morelos'
$ws.Range("B569").Value = 'Pinotepa De Don Luis'
$ws.Range("B571").Value = 'Putla Villa De Guerrero'
$ws.Range("B572").Value = 'Reforma De Pineda'
$ws.Range("B578").Value = 'San Agustín De Las Juntas'
$ws.Range("B587").Value = 'San Antonino El Alto'
$ws.Range("B592").Value = 'San Antonio De La Cal'
$ws.Range("B613").Value = 'San José Del Progreso'
$ws.Range("B618").Value = 'San Juan Bautista Lo De Soto'
$ws.Range("B652").Value = 'San Mateo Del Mar'
$ws.Range("B667").Value = 'San Miguel Del Puerto'
$ws.Range("B672").Value = 'San Pablo Villa De Mitla'
$ws.Range("B689").Value = 'San Pedro Y San Pablo Ayutla'
$ws.Range("B690").Value = 'San Pedro Y San Pablo Teposcolula'
$ws.Range("B691").Value = 'San Pedro Y San Pablo Tequixtepec'
$ws.Range("B707").Value = 'Santa Cruz Tacache De Mina'
$ws.Range("B712").Value = 'Santa Lucía Del Camino'
$ws.Range("B723").Value = 'Santa María Jalapa Del Marqués'
$ws.Range("B781").Value = 'Santo Domingo De Morelos'
$ws.Range("B789").Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range("B790").Value = 'Tataltepec De Valdés'
$ws.Range("B791").Value = 'Teococuilco De Marcos Pérez'
$ws.Range("B792").Value = 'Teotitlán De Flores Magón'
$ws.Range("B794").Value = 'Tezoatlán De Segura Y Luna'
$ws.Range("B795").Value = 'Tlacolula De Matamoros'
$ws.Range("B797").Value = 'Totontepec Villa De Morelos'
$ws.Range("B800").Value = 'Villa Sola De Vega'
$ws.Range("B801").Value = 'Villa Tejúpam De La Unión'
$ws.Range("B802").Value = 'Villa De Chilapa De Díaz'
$ws.Range("B803").Value = 'Villa De Etla'
$ws.Range("B804").Value = 'Villa De Tututepec'
$ws.Range("B805").Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range("B806").Value = 'Villa De Zaachila'
$ws.Range("B808").Value = 'Zimatlán De Álvarez'
$ws.Range("B832").Value = 'Chalchicomula De Sesma'
$ws.Range("B841").Value = 'Chila De La Sal'
$ws.Range("B853").Value = 'Cuayuca De Andrade'
$ws.Range("B854").Value = 'Cuetzalan Del Progreso'
$ws.Range("B871").Value = 'Huehuetlán El Chico'
$ws.Range("B872").Value = 'Huehuetlán El Grande'
$ws.Range("B876").Value = 'Ixcamilpa De Guerrero'
$ws.Range("B879").Value = 'Izúcar De Matamoros'
$ws.Range("B887").Value = 'Los Reyes De Juárez'
$ws.Range("B897").Value = 'Palmar De Bravo'
$ws.Range("B907").Value = 'San Diego La Mesa Tochimiltzingo'
$ws.Range("B920").Value = 'San Nicolás De Los Ranchos'
$ws.Range("B925").Value = 'San Salvador El Seco'
$ws.Range("B926").Value = 'San Salvador El Verde'
$ws.Range("B932").Value = 'Tecali De Herrera'
$ws.Range("B940").Value = 'Tepanco De López'
$ws.Range("B941").Value = 'Tepatlaxco De Hidalgo'
$ws.Range("B947").Value = 'Tepexi De Rodríguez'
$ws.Range("B949").Value = 'Tepeyahualco De Cuauhtémoc'
$ws.Range("B950").Value = 'Tetela De Ocampo'
$ws.Range("B955").Value = 'Tlacotepec De Benito Juárez'
$ws.Range("B970").Value = 'Xayacatlán De Bravo'
$ws.Range("B990").Value = 'Cadereyta De Montes'
$ws.Range("B991").Value = 'Jalpan De Serra'
$ws.Range("B992").Value = 'Landa De Matamoros'
$ws.Range("B993").Value = 'Pinal De Amoles'
$ws.Range("B995").Value = 'San Juan Del Río'
$ws.Range("B1014").Value = 'Villa De Guadalupe'
$ws.Range("B1041").Value = 'Jalpa De Méndez'
$ws.Range("B1054").Value = 'Soto La Marina'
$ws.Range("B1058").Value = 'Acuamanala De Miguel Hidalgo'
$ws.Range("B1059").Value = 'Amaxac De Guerrero'
$ws.Range("B1064").Value = 'Contla De Juan Cuamatzi'
$ws.Range("B1070").Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range("B1073").Value = 'Mazatecochco De José María Morelos'
$ws.Range("B1074").Value = 'Nanacamilpa De Mariano Arista'
$ws.Range("B1077").Value = 'Papalotla De Xicohténcatl'
$ws.Range("B1081").Value = 'San Pablo Del Monte'
$ws.Range("B1089").Value = 'Tepetitla De Lardizábal'
$ws.Range("B1092").Value = 'Tetla De La Solidaridad'
$ws.Range("A1104").Value = 'Veracruz De Ignacio De La Llave'
$ws.Range("B1108").Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range("B1117").Value = 'Boca Del Río'
$ws.Range("B1131").Value = 'Cosamaloapan De Carpio'
$ws.Range("B1148").Value = 'Hueyapan De Ocampo'
$ws.Range("B1149").Value = 'Huiloapan De Cuauhtémoc'
$ws.Range("B1150").Value = 'Ignacio De La Llave'
$ws.Range("B1154").Value = 'Ixhuacán De Los Reyes'
$ws.Range("B1155").Value = 'Ixhuatlán De Madero'
$ws.Range("B1156").Value = 'Ixhuatlán Del Café'
$ws.Range("B1157").Value = 'Ixhuatlán Del Sureste'
$ws.Range("B1164").Value = 'Juchique De Ferrer'
$ws.Range("B1169").Value = 'Martínez De La Torre'
$ws.Range("B1173").Value = 'Nanchital De Lázaro Cárdenas Del Río'
$ws.Range("B1180").Value = 'Paso De Ovejas'
$ws.Range("B1183").Value = 'Poza Rica De Hidalgo'
$ws.Range("B1190").Value = 'Sayula De Alemán'
$ws.Range("B1191").Value = 'Soledad De Doblado'
$ws.Range("B1196").Value = 'Tatahuicapan De Juárez'
$ws.Range("B1215").Value = 'Vega De Alatorre'
$ws.Range("B1223").Value = 'Zontecomatlán De López Y Fuentes'
$ws.Range("B1232").Value = 'Concepción Del Oro'
$ws.Range("B1235").Value = 'Nochistlán De Mejía'

# 3. Normalize the grand-total label
$ws.Range("A1241").Value = 'Total'

# 4. Remove the trailing footnote rows (1243-1247) that are no longer part of the clean dataset
$ws.Range("A1243:D1247").EntireRow.Delete()

